# Update the "Correspond Handoff Datetime" (D3) and
# "Correspond Handback DateTime" (G3) timestamps on the zh-cn and de-de
# report sheets, as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 15:27:54"
$wsZhCn.Range("G3").Value = "2016-01-08 15:28:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 15:28:08"
$wsDeDe.Range("G3").Value = "2016-01-08 15:29:07"
